# Merge the E/F score columns across each pair of rows belonging to the
# same "avaliador" sub-block (one row already holds the value, the other
# is its duplicate used only so the block renders as two visual rows).
# Merging collapses them into a single cell spanning both rows, clearing
# the now-redundant duplicate value, matching the existing merge pattern
# already used for the A:B and C:D columns in this sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowPairs = @(
    @(6, 7),
    @(8, 9),
    @(12, 13),
    @(14, 15),
    @(18, 19),
    @(20, 21),
    @(24, 25),
    @(26, 27),
    @(30, 31),
    @(32, 33),
    @(36, 37),
    @(38, 39),
    @(42, 43),
    @(44, 45),
    @(48, 49),
    @(50, 51),
    @(54, 55),
    @(56, 57)
)

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $ws.Range("E$r1" + ":E$r2").Merge() | Out-Null
    $ws.Range("F$r1" + ":F$r2").Merge() | Out-Null
}
